# SearchKeywords.xlsx — add a "Search Result" column next to the existing
# "Search Terms" column, giving each search term the URL it resolved to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + values for column B.
$ws.Range("B1").Value = "Search Result"
$ws.Range("B2").Value = "www.selenium.dev"
$ws.Range("B3").Value = "www.uipath.com"

# Touch the formatting of the whole used range (A1:B3) so it carries an
# explicit (non-default) cell style, same as the rest of the sheet.
$ws.Range("A1:B3").VerticalAlignment = -4107

# Give the new column a comfortable width to fit the URLs.
$ws.Columns("B:B").ColumnWidth = 19.21875

# Leave the selection where the author left it.
$ws.Range("G10").Select() | Out-Null
